# update create trans to v2.1
#
# The CHARs worksheet had its header row reorganised:
#   - the palaeo/lexical/morphsyn human-eval columns (6 headers) were
#     collapsed into 4 generic "he_human_N" headers
#   - the trailing two columns (old Z "line_status_end" / AA "commentary")
#     were removed, so every header from U onward shifts one "slot" left
#   - per-row data validation follows the same shift: the Hebrew-letter
#     list validation that used to cover Q:V now only covers Q:T, and the
#     DAMAGED/DAMAGED_STILL_READ/NOT_DAMAGED list validation that used to
#     cover X:Z now covers V:X.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CHARs")

# ---- 1. Re-label the header row (A1:Y1 keep their positions; only the
#         text of Q1..Y1 changes, and the old Z1/AA1 headers disappear). ----
$ws.Range("Q1").Value = "he_human_0"
$ws.Range("R1").Value = "he_human_1"
$ws.Range("S1").Value = "he_human_2"
$ws.Range("T1").Value = "he_human_3"
$ws.Range("U1").Value = "line_id"
$ws.Range("V1").Value = "line_status_int"
$ws.Range("W1").Value = "line_status_mid"
$ws.Range("X1").Value = "line_status_end"
$ws.Range("Y1").Value = "commentary"

# ---- 2. Drop the two trailing columns entirely (Z:AA), shrinking the
#         sheet's used range from A1:AA8 down to A1:Y8. ----
$ws.Range("Z1:AA8").Delete()

# ---- 3. Fix up the per-row data validation lists (rows 1-7). ----
for ($row = 1; $row -le 7; $row++) {
    # The Hebrew-letter picklist used to span Q:V; now it only spans Q:T,
    # so drop it from the U and V cells.
    $ws.Range("U" + $row).Validation.Delete()
    $ws.Range("V" + $row).Validation.Delete()

    # The DAMAGED/DAMAGED_STILL_READ/NOT_DAMAGED picklist used to span
    # X:Z; now it spans V:X. X already carries it, so drop the old Y
    # validation and add fresh ones on V and W.
    $ws.Range("Y" + $row).Validation.Delete()
    $ws.Range("V" + $row).Validation.Add(3, 1, 1, '"DAMAGED,DAMAGED_STILL_READ,NOT_DAMAGED"')
    $ws.Range("W" + $row).Validation.Add(3, 1, 1, '"DAMAGED,DAMAGED_STILL_READ,NOT_DAMAGED"')
}
